$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new checklist rows (same column/layout as the preceding
# "Sofia" entries at A19:A21).
$ws.Range("A22").Value = "Skapa ny användare - lägg till mejl"
$ws.Range("A23").Value = "När man skiftar mellan sidor sparas datan (fråga Sofia)"

# Match the formatting (fill/style) used by the rows directly above.
$ws.Range("A19").Copy()
$ws.Range("A22:A23").PasteSpecial(-4122)

# Move the active selection the way the author's session ended up.
$ws.Range("A27").Select()
